$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Update header row (A1:C1 renamed, new columns D1:I1 added) ---
$ws.Range("A1").Value = "Vendor"
$ws.Range("B1").Value = "FirstName"
$ws.Range("C1").Value = "LastName"
$ws.Range("D1").Value = "Item1"
$ws.Range("E1").Value = "Qty1"
$ws.Range("F1").Value = "Item2"
$ws.Range("G1").Value = "Qty2"
$ws.Range("H1").Value = "Item3"
$ws.Range("I1").Value = "Qty3"

# --- Vendor column (filled in this order to match original authoring order) ---
$ws.Range("A2").Value = "Amazon"
$ws.Range("A4").Value = "Target"
$ws.Range("A3").Value = "Walmart"
$ws.Range("A5").Value = "Amazon"

# --- Row 2: John / Doe / Burger / 2 ---
$ws.Range("B2").Value = "John "
$ws.Range("C2").Value = "Doe"
$ws.Range("D2").Value = "Burger"
$ws.Range("E2").Value = 2

# --- Row 3: Mary / Smith / Car / 1 ---
$ws.Range("B3").Value = "Mary "
$ws.Range("C3").Value = "Smith"
$ws.Range("D3").Value = "Car"
$ws.Range("E3").Value = 1

# --- Row 4: Tasha / Locke / Flowers / 5 ---
$ws.Range("B4").Value = "Tasha"
$ws.Range("C4").Value = "Locke"
$ws.Range("D4").Value = "Flowers"
$ws.Range("E4").Value = 5

# --- Row 5: Hassan / Baraka / Watch / 1 ---
$ws.Range("B5").Value = "Hassan"
$ws.Range("C5").Value = "Baraka"
$ws.Range("D5").Value = "Watch"
$ws.Range("E5").Value = 1

# --- Update selection to match the post-edit state ---
$ws.Range("E6").Select()
